# Uzupelnienie sprint_planu i InOUT
# Fill in "Dane wejsciowe" (D18) / "Dane wyjsciowe" (E18) for the HRV1 module
# row-pair (rows 18-19), merge each column across the two rows, and apply the
# header-style formatting (Verdana font for the "Dane wyjsciowe" column,
# centered + wrapped text, matching borders) that Excel applied to the rest
# of that two-row block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------------
$ws.Range("E18").Value = "map<String, dobule> - 7 parametrow analizy czasowej, map<String, dobule> - 6 parametrow analizy czestoliwosciowej, 2x vector<double> - oś X i oś Y wykresu oraz naniesione na niego paremetry z analizy czestoliwosciowej`n"
$ws.Range("D18").Value = "double - czestoliwosc probkowania, vector<double> - kolejne piki R z modulu R_peaks lub czas ich wystepienia"

# --- Formatting: "Dane wyjsciowe" column (E18:E19) --------------------------
$e18 = $ws.Range("E18")
$e18.Borders.Item(7).ColorIndex = 1
$e18.Borders.Item(7).LineStyle = 1
$e18.Borders.Item(7).Weight = 2
$e18.Borders.Item(8).ColorIndex = 1
$e18.Borders.Item(8).LineStyle = 1
$e18.Borders.Item(8).Weight = 2
$e18.Font.Name = "Verdana"
$e18.Font.Color = 0
$e18.HorizontalAlignment = -4108
$e18.VerticalAlignment = -4108
$e18.WrapText = $true

$e19 = $ws.Range("E19")
$e19.Borders.Item(7).ColorIndex = 1
$e19.Borders.Item(7).LineStyle = 1
$e19.Borders.Item(7).Weight = 2
$e19.Borders.Item(9).ColorIndex = 1
$e19.Borders.Item(9).LineStyle = 1
$e19.Borders.Item(9).Weight = 2
$e19.Font.Name = "Verdana"
$e19.Font.Color = 0
$e19.HorizontalAlignment = -4108
$e19.VerticalAlignment = -4108
$e19.WrapText = $true

# --- Formatting: "Dane wejsciowe" column (D18:D19) ---------------------------
$d18 = $ws.Range("D18")
$d18.HorizontalAlignment = -4108
$d18.WrapText = $true

# Row 18 picked up an automatic wrap-text height estimate from the edits
# above; re-running AutoFit settles it back on the sheet's default height
# (matching how Excel leaves a merged, wrapped row sized in the saved file).
$ws.Rows.Item(18).AutoFit()

$d19 = $ws.Range("D19")
$d19.HorizontalAlignment = -4108
$d19.WrapText = $true

# --- Merge the two-row blocks ----------------------------------------------
$ws.Range("E18:E19").Merge()
$ws.Range("D18:D19").Merge()

# --- Restore the cursor/selection position Excel left the sheet in ----------
$ws.Range("D22").Select() | Out-Null
